$d = $word.ActiveDocument

# Locate the target paragraph by a unique anchor substring (the "In Person" hearing-method
# paragraph that contains the repeated "further details will be provided" text).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text -like "*disposalHearingMethodInPerson*further details will be provided*") {
    $targetIndex = $i
    break
  }
}

if ($targetIndex -eq -1) {
  throw "Could not locate target paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)
$full = $p.Range
# Exclude the trailing paragraph mark so the paragraph properties (list numbering,
# spacing, etc.) are preserved and only the run content is replaced.
$r = $d.Range($full.Start, $full.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="0003704E"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;cs_{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0003704E"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>disposalHearingMethodDJ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0003704E"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>=’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00C62AEF"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>disposalHearingMethodInPerso</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>n</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="0003704E"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>’</w:t></w:r><w:r w:rsidR="00A15CEA"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>}&gt;&gt;</w:t></w:r><w:r w:rsidRPr="00430A35"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00430A35"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">This hearing will be in </w:t></w:r><w:r w:rsidRPr="00430A35"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>person</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;cs_{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>!=null}&gt;&gt;at &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.court_address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt; &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.postcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, &lt;&lt;es_&gt;&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>further details will be provided in your hearing notice</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">&lt;&lt;es_&gt;&gt;&lt;&lt;cs_{disposalHearingMethodDJ=’disposalHearingMethodTelephoneHearing’}&gt;&gt; </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">This hearing </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;cs_{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>!=null}&gt;&gt;will take place at &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.court_address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt; &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.postcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, and &lt;&lt;es_&gt;&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">will be attended by </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>telephone</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>, further details will be provided in your hearing notice</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;es_&gt;&gt;&lt;&lt;cs_{disposalHearingMethodDJ=’disposalHearingMethodVideoConferenceHearing’}&gt;&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">This hearing </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;cs_{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>!=null}&gt;&gt;will take place at &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.site_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.court_address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt; &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>hearingLocation.postcode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&gt;&gt;, and &lt;&lt;es_&gt;&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">will be attended by </w:t></w:r><w:r w:rsidRPr="0014274F"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>video conference</w:t></w:r><w:r w:rsidRPr="00430A35"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>, further details will be provided in your hearing notice</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidRPr="007D7981"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>&lt;&lt;es_&gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Host "Updated paragraph" $targetIndex
Write-Host $d.Paragraphs.Item($targetIndex).Range.Text
